# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps on the report sheets.

$wb = $excel.ActiveWorkbook

# "5b8f047b-...md" row: Overview!G3 and de-de!H3 both hold the same
# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" value
# (they shared the same underlying string), so both must be updated together.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-09-02 06:53:12"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-09-02 06:53:12"
$dede.Range("K3").Value = "2016-09-02 06:53:58"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-09-02 06:53:02"
$zhcn.Range("K3").Value = "2016-09-02 06:53:50"
